$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.347185373306274
$ws.Range("B1").Value = 2.156166553497314
$ws.Range("C1").Value = 4.883908271789551
$ws.Range("D1").Value = 3.268582820892334
$ws.Range("E1").Value = 1.286518812179565
